$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false,
                             $true, 1, $false, $replace, 2)
}

Replace-Text "2024-01-28 Sunday" "2024-01-29 Monday"

Replace-Text "50×14=700" "34×30=1020"
Replace-Text "59×55=3245" "99×30=2970"
Replace-Text "50×58=2900" "77×98=7546"
Replace-Text "84×62=5208" "46×31=1426"
Replace-Text "55×11=605" "23×87=2001"

Replace-Text "86×14=1204" "36×98=3528"
Replace-Text "85×91=7735" "11×20=220"
Replace-Text "73×29=2117" "67×51=3417"
Replace-Text "68×38=2584" "76×37=2812"
Replace-Text "98×96=9408" "36×59=2124"

Replace-Text "90×72=6480" "51×78=3978"
Replace-Text "97×29=2813" "48×70=3360"
Replace-Text "46×93=4278" "83×49=4067"
Replace-Text "50×63=3150" "96×16=1536"
Replace-Text "82×48=3936" "23×11=253"

Replace-Text "84×41=3444" "56×30=1680"
Replace-Text "28×45=1260" "22×86=1892"
Replace-Text "64×92=5888" "47×30=1410"
Replace-Text "18×50=900" "58×84=4872"
Replace-Text "82×97=7954" "58×99=5742"

Replace-Text "29×14=406" "87×61=5307"
Replace-Text "88×44=3872" "77×73=5621"
Replace-Text "71×19=1349" "61×55=3355"
Replace-Text "66×34=2244" "76×47=3572"
Replace-Text "76×86=6536" "28×96=2688"
